$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order of players (rows 2-19), with updated Position and Team values
$players = @(
    "Anfernee Simons",
    "Jayson Tatum",
    "Anthony Edwards",
    "Jaren Jackson Jr.",
    "Tyus Jones",
    "Rui Hachimura",
    "Aaron Gordon",
    "Ivica Zubac",
    "Giannis Antetokounmpo",
    "James Harden",
    "Collin Sexton",
    "Fred VanVleet",
    "Khris Middleton",
    "Paul George",
    "Jaden Ivey",
    "Zion Williamson",
    "Kyle Kuzma",
    "Andrew Wiggins"
)

$positions = @(
    "PG,SG",
    "SF,PF",
    "SG,SF",
    "PF,C",
    "PG",
    "SF,PF",
    "PF,C",
    "C",
    "PF,C",
    "PG,SG",
    "PG,SG",
    "PG",
    "SF",
    "SG,SF,PF",
    "PG,SG",
    "PF,C",
    "PF",
    "SF,PF"
)

$teams = @(
    "Portland Trail Blazers",
    "Boston Celtics",
    "Minnesota Timberwolves",
    "Memphis Grizzlies",
    "Phoenix Suns",
    "Los Angeles Lakers",
    "Denver Nuggets",
    "LA Clippers",
    "Milwaukee Bucks",
    "LA Clippers",
    "Utah Jazz",
    "Houston Rockets",
    "Milwaukee Bucks",
    "Philadelphia 76ers",
    "Detroit Pistons",
    "New Orleans Pelicans",
    "Washington Wizards",
    "Golden State Warriors"
)

for ($i = 0; $i -lt $players.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $players[$i]
}

for ($i = 0; $i -lt $positions.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $positions[$i]
}

for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $teams[$i]
}
